$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "MN"
$ws.Range("H2").Value = "TN"
$ws.Range("R3").Value = "KN"
$ws.Range("G14").Value = "PN"
$ws.Range("B14").Value = "CN"
$ws.Range("B20").Value = "ON"
$ws.Range("L20").Value = "UN"
$ws.Range("H7").Value = "DN"
$ws.Range("R13").Value = "LN"
$ws.Range("B8").Value = "VN"

$ws.Range("B8").Select()
